# Update master to output generated at 9a8706d
$d = $word.ActiveDocument

$replacements = @(
    @("442×4=1768", "445×4=1780"),
    @("959×9=8631", "395×6=2370"),
    @("240×6=1440", "239×3=717"),
    @("476×5=2380", "942×5=4710"),
    @("303×8=2424", "433×7=3031"),
    @("623×9=5607", "203×6=1218"),
    @("320×8=2560", "681×8=5448"),
    @("606×2=1212", "248×7=1736"),
    @("374×7=2618", "998×7=6986"),
    @("638×7=4466", "802×6=4812"),
    @("709×9=6381", "776×9=6984"),
    @("375×2=750", "234×6=1404"),
    @("737×4=2948", "177×4=708"),
    @("468×9=4212", "796×5=3980"),
    @("233×4=932", "990×5=4950"),
    @("548×4=2192", "249×5=1245"),
    @("606×6=3636", "120×8=960"),
    @("643×7=4501", "176×7=1232"),
    @("367×9=3303", "867×8=6936"),
    @("797×7=5579", "331×4=1324"),
    @("632×8=5056", "851×9=7659"),
    @("998×2=1996", "604×7=4228"),
    @("406×9=3654", "224×4=896"),
    @("874×9=7866", "494×4=1976"),
    @("310×6=1860", "460×4=1840"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
